$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.496.52"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.918.78"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("D7").Value = "'0.4803"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "'0.08210"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'23.44"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.921.25"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "'6.052"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "'7.229"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "'91.36"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "'0.06845"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "'0.00001038"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'17.54"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "29.498.14"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("D23").Value = "'11.88"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "2.156.50"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'156.28"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").Value = "'6.494"
$ws.Range("E27").Value = "  +3.09%  "
$ws.Range("D28").Value = "'19.96"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").Value = "'2.095"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "'120.47"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").Value = "'0.09610"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "'5.611"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").Value = "'3.558"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "'1.370"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").Value = "'0.06324"
$ws.Range("E36").Value = "  +3.53%  "
$ws.Range("D37").Value = "'0.02282"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").Value = "'1.183"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "'0.5926"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "'10.72"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'7.882"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").Value = "'0.1846"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "'2.392"
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "'0.07474"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").Value = "'0.5560"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'1.934"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "'117.98"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").Value = "'2.427"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("D51").Value = "'71.98"
$ws.Range("E51").Value = "  -0.76%  "
